$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12): correct the scoring numbers, which were
# previously wrong because of duplicate question columns (see below). ---

# Row 10 ("No."): Right=19, Wrong=0 (unchanged), Not Attempt=9, Max=28
$ws.Range("B10").Value = 19
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 28

# Row 11 ("Marking"): Right mark=4, Wrong mark=-1 as a real number (was text "-1")
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 ("Total"): Right total=76, and the Max-string becomes "76/112"
$ws.Range("B12").Value = 76
$ws.Range("E12").Value = "76/112"

# Give A10, A11, A12 the same ("correctStyle"/green, s=4) formatting already
# used by the row-9 header cells, instead of the default formatting.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

# --- Remove the duplicated 3rd "Student Ans / Correct Ans" block (columns
# G:H) entirely -- it was a full duplicate of the 1st/2nd block. Deleting
# the whole columns also shrinks the sheet dimension from A5:H40 to A5:E40. ---
$ws.Range("F1:H40").EntireColumn.Delete()

# --- Remove the duplicated 2nd "Student Ans / Correct Ans" block (columns
# D:E) for every question row except the first two, which stay. ---
$ws.Range("D19:E40").Clear()

# Fill in D16 and D17 (2nd block "Student Ans") with the text that matches
# the "Correct Ans" in E16/E17, using the same red "incorrectStyle" (s=5)
# formatting already used by the 1st block's filled-in Student-Ans cells.
$ws.Range("B10").Copy()
$ws.Range("D16:D17").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"

# --- Fill in column A (1st block "Student Ans") for every question row
# where the student's recorded answer equals the correct answer, copying
# the formatting (s=5) the same way. Rows left blank keep their original
# empty "normalStyle" (s=7) placeholder cell untouched. ---
$answers = @{
    16 = "Option A"
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    22 = "Option D"
    23 = "Option D"
    28 = "Option D"
    29 = "Option D"
    30 = "Option B"
    31 = "Option D"
    32 = "Option C"
    33 = "Option D"
    35 = "Option D"
    36 = "Option A"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
}
foreach ($r in $answers.Keys) {
    $ws.Range("B10").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = $answers[$r]
}
